$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 71666.44500000001
$ws.Range("J87").Value = 71666.44500000001
$ws.Range("L87").Value = 71666.44500000001
$ws.Range("N87").Value = -74162.44500000001
$ws.Range("H90").Value = 71666.44500000001
$ws.Range("J90").Value = 71666.44500000001
$ws.Range("L90").Value = 214999.335
$ws.Range("N90").Value = -227479.335
$ws.Range("H109").Value = 49999.5
$ws.Range("J109").Value = 49999.5
$ws.Range("L109").Value = 49999.5
$ws.Range("N109").Value = -52773.5
$ws.Range("H126").Value = 50799
$ws.Range("J126").Value = 50799
$ws.Range("L126").Value = 50799
$ws.Range("N126").Value = -60679
$ws.Range("H130").Value = 79990
$ws.Range("J130").Value = 79990
$ws.Range("L130").Value = 79990
$ws.Range("N130").Value = -90030
$ws.Range("H132").Value = 1465.5128
$ws.Range("I132").Value = 1364.0278
$ws.Range("K132").Value = 4092.0834
$ws.Range("M132").Value = -1562.0834
$ws.Range("H138").Value = 5356.2188
$ws.Range("I138").Value = 1848.4166
$ws.Range("J138").Value = 6165.7114
$ws.Range("K138").Value = 5545.2498
$ws.Range("L138").Value = 18497.1342
$ws.Range("M138").Value = -405.2497999999996
$ws.Range("N138").Value = -28777.1342
$ws.Range("H139").Value = 96494.836
$ws.Range("J139").Value = 99994
$ws.Range("L139").Value = 99994
$ws.Range("N139").Value = -110274
$ws.Range("H141").Value = 6411730.5
$ws.Range("J141").Value = 3595
$ws.Range("L141").Value = 10785
$ws.Range("N141").Value = -21145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 877.625
$ws.Range("K2").Value = 877.625
$ws.Range("M2").Value = -764.625
$ws.Range("H61").Value = 9303.723
$ws.Range("I61").Value = 859.7143
$ws.Range("J61").Value = 14677.182
$ws.Range("K61").Value = 859.7143
$ws.Range("L61").Value = 14677.182
$ws.Range("M61").Value = -647.7143
$ws.Range("N61").Value = -15101.182
$ws.Range("H63").Value = 1332.5714
$ws.Range("I63").Value = 1107
$ws.Range("K63").Value = 1107
$ws.Range("M63").Value = -421
$ws.Range("H66").Value = 1332.5714
$ws.Range("I66").Value = 1107
$ws.Range("K66").Value = 5535
$ws.Range("M66").Value = -2103
$ws.Range("H105").Value = 400000
$ws.Range("J105").Value = 400000
$ws.Range("L105").Value = 400000
$ws.Range("N105").Value = -406988
$ws.Range("I116").Value = 877.625
$ws.Range("K116").Value = 877.625
$ws.Range("M116").Value = 1416.375
$ws.Range("H122").Value = 2386.9
$ws.Range("I122").Value = 1466.963
$ws.Range("J122").Value = 10666.333
$ws.Range("K122").Value = 4400.889
$ws.Range("L122").Value = 31998.999
$ws.Range("M122").Value = -1950.889
$ws.Range("N122").Value = -36898.999
$ws.Range("H132").Value = 5356.893
$ws.Range("I132").Value = 4407.6665
$ws.Range("J132").Value = 6452.154
$ws.Range("K132").Value = 13222.9995
$ws.Range("L132").Value = 19356.462
$ws.Range("M132").Value = -10692.9995
$ws.Range("N132").Value = -24416.462
$ws.Range("H136").Value = 9303.723
$ws.Range("I136").Value = 859.7143
$ws.Range("J136").Value = 14677.182
$ws.Range("K136").Value = 2579.1429
$ws.Range("L136").Value = 44031.546
$ws.Range("M136").Value = -29.14289999999983
$ws.Range("N136").Value = -49131.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 877.625
$ws.Range("K3").Value = 877.625
$ws.Range("M3").Value = -763.625
$ws.Range("H20").Value = 8773966
$ws.Range("I20").Value = 13891018
$ws.Range("J20").Value = 1878
$ws.Range("K20").Value = 13891018
$ws.Range("L20").Value = 1878
$ws.Range("M20").Value = -13890771
$ws.Range("N20").Value = -2372
$ws.Range("H99").Value = 11366162
$ws.Range("I99").Value = 2100
$ws.Range("J99").Value = 22730224
$ws.Range("K99").Value = 2100
$ws.Range("L99").Value = 22730224
$ws.Range("M99").Value = -602
$ws.Range("N99").Value = -22733220
$ws.Range("H134").Value = 5787.9165
$ws.Range("I134").Value = 2086.5
$ws.Range("K134").Value = 6259.5
$ws.Range("M134").Value = -3724.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6438.1763
$ws.Range("I16").Value = 5459
$ws.Range("K16").Value = 5459
$ws.Range("M16").Value = -5172
$ws.Range("H31").Value = 5400.62
$ws.Range("I31").Value = 2518.4038
$ws.Range("K31").Value = 2518.4038
$ws.Range("M31").Value = -2223.4038
$ws.Range("H34").Value = 5400.62
$ws.Range("I34").Value = 2518.4038
$ws.Range("K34").Value = 2518.4038
$ws.Range("M34").Value = -2316.4038
$ws.Range("H113").Value = 6438.1763
$ws.Range("I113").Value = 5459
$ws.Range("K113").Value = 5459
$ws.Range("M113").Value = -3289
$ws.Range("H132").Value = 5472.104
$ws.Range("I132").Value = 2896.8928
$ws.Range("K132").Value = 8690.678400000001
$ws.Range("M132").Value = -6160.678400000001
$ws.Range("H134").Value = 4392.4326
$ws.Range("I134").Value = 1961.0476
$ws.Range("J134").Value = 8477.16
$ws.Range("K134").Value = 5883.142800000001
$ws.Range("L134").Value = 25431.48
$ws.Range("M134").Value = -3348.142800000001
$ws.Range("N134").Value = -30501.48

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3333724
$ws.Range("I12").Value = 905
$ws.Range("J12").Value = 5000133.5
$ws.Range("K12").Value = 2715
$ws.Range("L12").Value = 15000400.5
$ws.Range("M12").Value = -2542
$ws.Range("N12").Value = -15000746.5
$ws.Range("H37").Value = 93999.664
$ws.Range("J37").Value = 93999.664
$ws.Range("L37").Value = 281998.992
$ws.Range("N37").Value = -282222.992
$ws.Range("H76").Value = 5943.5
$ws.Range("J76").Value = 8888
$ws.Range("L76").Value = 26664
$ws.Range("N76").Value = -27430
$ws.Range("H79").Value = 5943.5
$ws.Range("J79").Value = 8888
$ws.Range("L79").Value = 26664
$ws.Range("N79").Value = -29316
$ws.Range("H80").Value = 40004080
$ws.Range("I80").Value = 27781612
$ws.Range("K80").Value = 83344836
$ws.Range("M80").Value = -83343900
$ws.Range("H83").Value = 40004080
$ws.Range("I83").Value = 27781612
$ws.Range("K83").Value = 250034508
$ws.Range("M83").Value = -250029828
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1984.75
$ws.Range("I80").Value = 1879.6666
$ws.Range("K80").Value = 1879.6666
$ws.Range("M80").Value = -881.6666
$ws.Range("H83").Value = 1984.75
$ws.Range("I83").Value = 1879.6666
$ws.Range("K83").Value = 9398.333000000001
$ws.Range("M83").Value = -4406.333000000001
$ws.Range("H107").Value = 1333589.6
$ws.Range("I107").Value = 1600212.6
$ws.Range("K107").Value = 1600212.6
$ws.Range("M107").Value = -1598292.6
$ws.Range("H113").Value = 328989.53
$ws.Range("I113").Value = 1002490
$ws.Range("J113").Value = 8275
$ws.Range("K113").Value = 1002490
$ws.Range("L113").Value = 8275
$ws.Range("M113").Value = -1000320
$ws.Range("N113").Value = -12615
$ws.Range("H132").Value = 5406.8105
$ws.Range("I132").Value = 2418.9714
$ws.Range("J132").Value = 9953.521000000001
$ws.Range("K132").Value = 7256.914199999999
$ws.Range("L132").Value = 29860.563
$ws.Range("M132").Value = -4726.914199999999
$ws.Range("N132").Value = -34920.563

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4594.737
$ws.Range("I93").Value = 2554.5454
$ws.Range("J93").Value = 7400
$ws.Range("K93").Value = 2554.5454
$ws.Range("L93").Value = 7400
$ws.Range("M93").Value = -1306.5454
$ws.Range("N93").Value = -9896
$ws.Range("H100").Value = 4645.8
$ws.Range("I100").Value = 3314.6667
$ws.Range("K100").Value = 3314.6667
$ws.Range("M100").Value = -2773.6667
$ws.Range("H122").Value = 4787.943
$ws.Range("I122").Value = 3155.348
$ws.Range("J122").Value = 7917.0835
$ws.Range("K122").Value = 9466.044
$ws.Range("L122").Value = 23751.2505
$ws.Range("M122").Value = -7016.044
$ws.Range("N122").Value = -28651.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 936.2308
$ws.Range("I100").Value = 550
$ws.Range("J100").Value = 1107.8889
$ws.Range("K100").Value = 1100
$ws.Range("L100").Value = 2215.7778
$ws.Range("M100").Value = -559
$ws.Range("N100").Value = -3297.7778
$ws.Range("H113").Value = 12150.091
$ws.Range("I113").Value = 14716.833
$ws.Range("K113").Value = 44150.499
$ws.Range("M113").Value = -41980.499
$ws.Range("H122").Value = 157389.23
$ws.Range("I122").Value = 268213.66
$ws.Range("K122").Value = 804640.98
$ws.Range("M122").Value = -802190.98
$ws.Range("H132").Value = 15165060
$ws.Range("I132").Value = 25007148
$ws.Range("J132").Value = 23388.23
$ws.Range("K132").Value = 75021444
$ws.Range("L132").Value = 70164.69
$ws.Range("M132").Value = -75018914
$ws.Range("N132").Value = -75224.69
$ws.Range("H136").Value = 32295436
$ws.Range("I136").Value = 71429816
$ws.Range("J136").Value = 67122.88
$ws.Range("K136").Value = 214289448
$ws.Range("L136").Value = 201368.64
$ws.Range("M136").Value = -214286898
$ws.Range("N136").Value = -206468.64
